$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B7").Value = 16052082.46800645
$ws.Range("B8").Value = 16052082.46800645
$ws.Range("B9").Value = 3940022.790763975
$ws.Range("B10").Value = 3940022.790763975
$ws.Range("B12").Value = 59583988.52783271
$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("H2").Value = 347.8590406130752
$ws.Range("S2").Value = 359.6204767756762
$ws.Range("G3").Value = 322.5970075731406
$ws.Range("I3").Value = 107.3959878205679
$ws.Range("S5").Value = 361.0088210686545
$ws.Range("X5").Value = 590.8934891676897
$ws.Range("G6").Value = 322.5970075731406
$ws.Range("I6").Value = 107.3959878205679
$ws.Range("H8").Value = 347.8590406130752
$ws.Range("S8").Value = 359.6204767756762
$ws.Range("F9").Value = 338.1025423803039
$ws.Range("I9").Value = 108.929687778141
$ws.Range("R10").Value = 208.1413738827904
$ws.Range("F11").Value = 403.5012844150336
$ws.Range("Y11").Value = 511.3174326828064
$ws.Range("R12").Value = 352.1613821980396
$ws.Range("S12").Value = 406.6592151013782
$ws.Range("T13").Value = 165.2127885434264
$ws.Range("W13").Value = 0
$ws.Range("G14").Value = 396.4956611598672
$ws.Range("T14").Value = 560.964918411981
$ws.Range("R15").Value = 350.6276822404662
$ws.Range("S15").Value = 408.192915058951
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 165.2127885434264
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 165.2127885434264
$ws.Range("H20").Value = 346.4706963200988
$ws.Range("S20").Value = 361.0088210686545
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0
$ws.Range("R22").Value = 0
$ws.Range("S22").Value = 0
$ws.Range("T22").Value = 165.2127885434264
$ws.Range("R24").Value = 350.6276822404662
$ws.Range("S24").Value = 408.192915058951
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = 0
$ws.Range("S25").Value = 0
$ws.Range("T25").Value = 165.2127885434264
$ws.Range("H26").Value = 347.8590406130752
$ws.Range("X26").Value = 590.8934891676897
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("V28").Value = 165.2127885434264
$ws.Range("H29").Value = 347.8590406130752
$ws.Range("V29").Value = 628.4626798738458
$ws.Range("H31").Value = 165.2127885434264
$ws.Range("I31").Value = 0
$ws.Range("B32").Value = 481.9993129555745
$ws.Range("S32").Value = 359.6204767756762
$ws.Range("H34").Value = 165.2127885434264
$ws.Range("Y34").Value = 0
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 165.2127885434264
$ws.Range("G38").Value = 396.4956611598672
$ws.Range("Y38").Value = 511.3174326828064
$ws.Range("R39").Value = 350.6276822404662
$ws.Range("W39").Value = 432.3731429098285
$ws.Range("V40").Value = 165.2127885434264
$ws.Range("X40").Value = 0
$ws.Range("C41").Value = 449.4745782429939
$ws.Range("Y41").Value = 509.9290883898282
$ws.Range("R42").Value = 352.1613821980396
$ws.Range("S42").Value = 406.6592151013782
$ws.Range("G44").Value = 397.8840054528454
$ws.Range("Y44").Value = 509.9290883898282
$ws.Range("I46").Value = 97.40414414470841
$ws.Range("Q46").Value = 0
$ws.Range("T46").Value = 67.808644398718
$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("B2").Value = 2580.161715518833
$ws.Range("C2").Value = 2126.146990020859
$ws.Range("D2").Value = 1711.662994324041
$ws.Range("E2").Value = 1303.215227044376
$ws.Range("F2").Value = 894.2358041069906
$ws.Range("G2").Value = 492.3327682960356
$ws.Range("I2").Value = 140.96
$ws.Range("J2").Value = 599.7771454756521
$ws.Range("K2").Value = 1339.420289043491
$ws.Range("L2").Value = 2257.013012360075
$ws.Range("M2").Value = 2739.192211873577
$ws.Range("N2").Value = 3304.235275746249
$ws.Range("O2").Value = 4776.556370989979
$ws.Range("P2").Value = 6474.295852434775
$ws.Range("S2").Value = 6684.746993155883
$ws.Range("T2").Value = 6118.11576243671
$ws.Range("U2").Value = 5462.811290515621
$ws.Range("V2").Value = 4826.598134791557
$ws.Range("W2").Value = 4181.776441881509
$ws.Range("X2").Value = 3583.51196363841
$ws.Range("Y2").Value = 3067.029708403252
$ws.Range("G3").Value = 554.3073408685043
$ws.Range("H3").Value = 249.4407957783515
$ws.Range("I4").Value = 140.96
$ws.Range("J4").Value = 290.4466099967136
$ws.Range("K4").Value = 307.8416045893196
$ws.Range("L4").Value = 307.8416045893196
$ws.Range("I5").Value = 140.96
$ws.Range("J5").Value = 599.7771454756521
$ws.Range("K5").Value = 2049.12665539136
$ws.Range("L5").Value = 3681.051378010443
$ws.Range("M5").Value = 4163.230577523946
$ws.Range("N5").Value = 4728.273641396618
$ws.Range("O5").Value = 5638.428524337362
$ws.Range("P5").Value = 6474.295852434775
$ws.Range("S5").Value = 6683.344625183177
$ws.Range("T5").Value = 6116.713394464005
$ws.Range("U5").Value = 5461.408922542916
$ws.Range("V5").Value = 4825.195766818852
$ws.Range("W5").Value = 4180.374073908803
$ws.Range("G6").Value = 554.3073408685043
$ws.Range("H6").Value = 249.4407957783515
$ws.Range("B7").Value = 259.7315161073162
$ws.Range("C7").Value = 259.7315161073162
$ws.Range("D7").Value = 259.7315161073162
$ws.Range("E7").Value = 259.7315161073162
$ws.Range("F7").Value = 259.7315161073162
$ws.Range("G7").Value = 259.7315161073162
$ws.Range("H7").Value = 259.7315161073162
$ws.Range("I7").Value = 259.7315161073162
$ws.Range("J7").Value = 409.2181261040298
$ws.Range("K7").Value = 426.6131206966359
$ws.Range("L7").Value = 426.6131206966359
$ws.Range("M7").Value = 426.6131206966359
$ws.Range("N7").Value = 426.6131206966359
$ws.Range("S7").Value = 223.6648391925338
$ws.Range("T7").Value = 259.7315161073162
$ws.Range("U7").Value = 259.7315161073162
$ws.Range("V7").Value = 259.7315161073162
$ws.Range("W7").Value = 259.7315161073162
$ws.Range("X7").Value = 259.7315161073162
$ws.Range("Y7").Value = 259.7315161073162
$ws.Range("B8").Value = 2580.161715518833
$ws.Range("C8").Value = 2126.146990020859
$ws.Range("D8").Value = 1711.662994324041
$ws.Range("E8").Value = 1303.215227044376
$ws.Range("F8").Value = 894.2358041069906
$ws.Range("G8").Value = 492.3327682960356
$ws.Range("J8").Value = 675.3654425154225
$ws.Range("K8").Value = 2124.71495243113
$ws.Range("L8").Value = 3681.051378010443
$ws.Range("M8").Value = 4163.230577523946
$ws.Range("N8").Value = 4728.273641396618
$ws.Range("O8").Value = 5638.428524337362
$ws.Range("P8").Value = 6474.295852434775
$ws.Range("S8").Value = 6684.746993155883
$ws.Range("T8").Value = 6118.11576243671
$ws.Range("U8").Value = 5462.811290515621
$ws.Range("V8").Value = 4826.598134791557
$ws.Range("W8").Value = 4181.776441881509
$ws.Range("X8").Value = 3583.51196363841
$ws.Range("Y8").Value = 3067.029708403252
$ws.Range("F9").Value = 881.7120959500332
$ws.Range("G9").Value = 555.8565327448407
$ws.Range("H9").Value = 250.9899876546879
$ws.Range("K10").Value = 457.4949902104934
$ws.Range("L10").Value = 457.4949902104934
$ws.Range("M10").Value = 457.4949902104934
$ws.Range("N10").Value = 457.4949902104934
$ws.Range("O10").Value = 457.4949902104934
$ws.Range("P10").Value = 457.4949902104934
$ws.Range("B11").Value = 2578.759347546128
$ws.Range("C11").Value = 2124.744622048154
$ws.Range("D11").Value = 1710.260626351336
$ws.Range("E11").Value = 1301.812859071671
$ws.Range("I11").Value = 216.5482970397704
$ws.Range("J11").Value = 675.3654425154225
$ws.Range("K11").Value = 2124.71495243113
$ws.Range("L11").Value = 3071.357825473565
$ws.Range("M11").Value = 3553.537024987068
$ws.Range("N11").Value = 4118.58008885974
$ws.Range("O11").Value = 5028.734971800483
$ws.Range("P11").Value = 5864.602299897897
$ws.Range("Q11").Value = 6931.116362236997
$ws.Range("Y11").Value = 3065.627340430547
$ws.Range("R12").Value = 5569.989306372909
$ws.Range("T13").Value = 140.96
$ws.Range("U13").Value = 140.96
$ws.Range("V13").Value = 140.96
$ws.Range("B14").Value = 2578.759347546128
$ws.Range("C14").Value = 2124.744622048154
$ws.Range("D14").Value = 1710.260626351336
$ws.Range("E14").Value = 1301.812859071671
$ws.Range("F14").Value = 892.8334361342853
$ws.Range("I14").Value = 140.96
$ws.Range("J14").Value = 599.7771454756521
$ws.Range("K14").Value = 1830.673957362034
$ws.Range("L14").Value = 2748.266680678617
$ws.Range("M14").Value = 4280.114215286949
$ws.Range("N14").Value = 4845.157279159621
$ws.Range("O14").Value = 5755.312162100365
$ws.Range("P14").Value = 6591.179490197778
$ws.Range("Q14").Value = 7048
$ws.Range("T14").Value = 6116.713394464005
$ws.Range("U14").Value = 5461.408922542916
$ws.Range("V14").Value = 4825.195766818852
$ws.Range("W14").Value = 4180.374073908803
$ws.Range("X14").Value = 3582.109595665705
$ws.Range("Y14").Value = 3065.627340430547
$ws.Range("R15").Value = 5571.538498249245
$ws.Range("F16").Value = 307.8416045893196
$ws.Range("G16").Value = 307.8416045893196
$ws.Range("N17").Value = 4389.93206401015
$ws.Range("O17").Value = 5638.428524337362
$ws.Range("D19").Value = 307.8416045893196
$ws.Range("E19").Value = 307.8416045893196
$ws.Range("F19").Value = 307.8416045893196
$ws.Range("G19").Value = 307.8416045893196
$ws.Range("B20").Value = 2578.75934754613
$ws.Range("C20").Value = 2124.744622048156
$ws.Range("D20").Value = 1710.260626351338
$ws.Range("E20").Value = 1301.812859071673
$ws.Range("F20").Value = 892.8334361342871
$ws.Range("G20").Value = 490.9304003233322
$ws.Range("I20").Value = 216.5482970397704
$ws.Range("J20").Value = 1101.367131735718
$ws.Range("K20").Value = 2153.765102156984
$ws.Range("L20").Value = 3071.357825473567
$ws.Range("M20").Value = 3553.537024987069
$ws.Range("N20").Value = 4118.580088859741
$ws.Range("O20").Value = 5028.734971800485
$ws.Range("P20").Value = 5864.602299897899
$ws.Range("Q20").Value = 6931.116362236999
$ws.Range("S20").Value = 6683.344625183179
$ws.Range("T20").Value = 6116.713394464006
$ws.Range("U20").Value = 5461.408922542918
$ws.Range("V20").Value = 4825.195766818853
$ws.Range("W20").Value = 4180.374073908805
$ws.Range("X20").Value = 3582.109595665706
$ws.Range("Y20").Value = 3065.627340430548
$ws.Range("B22").Value = 140.96
$ws.Range("C22").Value = 140.96
$ws.Range("D22").Value = 140.96
$ws.Range("E22").Value = 140.96
$ws.Range("F22").Value = 140.96
$ws.Range("G22").Value = 140.96
$ws.Range("T22").Value = 140.96
$ws.Range("U22").Value = 140.96
$ws.Range("V22").Value = 140.96
$ws.Range("W22").Value = 140.96
$ws.Range("X22").Value = 140.96
$ws.Range("Y22").Value = 140.96
$ws.Range("I23").Value = 140.96
$ws.Range("J23").Value = 599.7771454756521
$ws.Range("K23").Value = 1339.420289043491
$ws.Range("L23").Value = 2257.013012360075
$ws.Range("M23").Value = 3553.537024987068
$ws.Range("N23").Value = 4118.58008885974
$ws.Range("O23").Value = 5028.734971800483
$ws.Range("P23").Value = 5864.602299897897
$ws.Range("R24").Value = 5571.538498249245
$ws.Range("Q25").Value = 307.8416045893196
$ws.Range("R25").Value = 307.8416045893196
$ws.Range("S25").Value = 307.8416045893196
$ws.Range("B26").Value = 2580.161715518833
$ws.Range("C26").Value = 2126.146990020859
$ws.Range("D26").Value = 1711.662994324041
$ws.Range("E26").Value = 1303.215227044376
$ws.Range("F26").Value = 894.2358041069906
$ws.Range("G26").Value = 492.3327682960356
$ws.Range("X26").Value = 3583.51196363841
$ws.Range("Y26").Value = 3067.029708403252
$ws.Range("B28").Value = 140.96
$ws.Range("C28").Value = 140.96
$ws.Range("D28").Value = 140.96
$ws.Range("E28").Value = 140.96
$ws.Range("F28").Value = 140.96
$ws.Range("G28").Value = 140.96
$ws.Range("H28").Value = 140.96
$ws.Range("V28").Value = 140.96
$ws.Range("W28").Value = 140.96
$ws.Range("X28").Value = 140.96
$ws.Range("Y28").Value = 140.96
$ws.Range("B29").Value = 2580.161715518833
$ws.Range("C29").Value = 2126.146990020859
$ws.Range("D29").Value = 1711.662994324041
$ws.Range("E29").Value = 1303.215227044376
$ws.Range("F29").Value = 894.2358041069906
$ws.Range("G29").Value = 492.3327682960356
$ws.Range("J29").Value = 675.3654425154225
$ws.Range("K29").Value = 1415.008586083262
$ws.Range("L29").Value = 2332.601309399844
$ws.Range("M29").Value = 3077.533789260046
$ws.Range("V29").Value = 4826.598134791557
$ws.Range("W29").Value = 4181.776441881509
$ws.Range("X29").Value = 3583.51196363841
$ws.Range("Y29").Value = 3067.029708403252
$ws.Range("H31").Value = 140.96
$ws.Range("I32").Value = 216.5482970397704
$ws.Range("J32").Value = 675.3654425154225
$ws.Range("K32").Value = 1713.790319599031
$ws.Range("L32").Value = 2631.383042915614
$ws.Range("M32").Value = 4163.230577523946
$ws.Range("N32").Value = 4728.273641396618
$ws.Range("O32").Value = 5638.428524337362
$ws.Range("S32").Value = 6684.746993155883
$ws.Range("T32").Value = 6118.11576243671
$ws.Range("U32").Value = 5462.811290515621
$ws.Range("V32").Value = 4826.598134791557
$ws.Range("W32").Value = 4181.776441881509
$ws.Range("X32").Value = 3583.51196363841
$ws.Range("Y32").Value = 3067.029708403252
$ws.Range("B34").Value = 307.8416045893196
$ws.Range("C34").Value = 307.8416045893196
$ws.Range("D34").Value = 307.8416045893196
$ws.Range("E34").Value = 307.8416045893196
$ws.Range("F34").Value = 307.8416045893196
$ws.Range("G34").Value = 307.8416045893196
$ws.Range("Y34").Value = 307.8416045893196
$ws.Range("J35").Value = 675.3654425154225
$ws.Range("K35").Value = 1415.008586083262
$ws.Range("L35").Value = 2332.601309399844
$ws.Range("M35").Value = 3077.533789260046
$ws.Range("C37").Value = 307.8416045893196
$ws.Range("D37").Value = 307.8416045893196
$ws.Range("E37").Value = 307.8416045893196
$ws.Range("F37").Value = 307.8416045893196
$ws.Range("G37").Value = 307.8416045893196
$ws.Range("B38").Value = 2578.759347546128
$ws.Range("C38").Value = 2124.744622048154
$ws.Range("D38").Value = 1710.260626351336
$ws.Range("E38").Value = 1301.812859071671
$ws.Range("F38").Value = 892.8334361342853
$ws.Range("I38").Value = 216.5482970397704
$ws.Range("J38").Value = 675.3654425154225
$ws.Range("K38").Value = 2018.47013910948
$ws.Range("L38").Value = 2936.062862426063
$ws.Range("M38").Value = 3418.242061939566
$ws.Range("N38").Value = 3983.285125812238
$ws.Range("O38").Value = 4893.440008752982
$ws.Range("P38").Value = 6591.179490197778
$ws.Range("Q38").Value = 7048
$ws.Range("Y38").Value = 3065.627340430547
$ws.Range("R39").Value = 5571.538498249245
$ws.Range("S39").Value = 5159.222422432123
$ws.Range("T39").Value = 4762.564674070616
$ws.Range("U39").Value = 4358.520619043174
$ws.Range("V39").Value = 3939.822975415375
$ws.Range("V40").Value = 140.96
$ws.Range("W40").Value = 140.96
$ws.Range("B41").Value = 2580.161715518833
$ws.Range("Y41").Value = 3067.029708403252
$ws.Range("R42").Value = 5569.989306372909
$ws.Range("B44").Value = 2580.161715518833
$ws.Range("C44").Value = 2126.146990020859
$ws.Range("D44").Value = 1711.662994324041
$ws.Range("E44").Value = 1303.215227044376
$ws.Range("F44").Value = 894.2358041069906
$ws.Range("I44").Value = 216.5482970397704
$ws.Range("J44").Value = 1101.367131735718
$ws.Range("K44").Value = 2153.765102156982
$ws.Range("L44").Value = 3071.357825473565
$ws.Range("M44").Value = 3553.537024987068
$ws.Range("N44").Value = 4118.58008885974
$ws.Range("O44").Value = 5028.734971800483
$ws.Range("P44").Value = 5864.602299897897
$ws.Range("Y44").Value = 3067.029708403252
$ws.Range("B46").Value = 239.3480243885944
$ws.Range("C46").Value = 239.3480243885944
$ws.Range("D46").Value = 239.3480243885944
$ws.Range("E46").Value = 239.3480243885944
$ws.Range("F46").Value = 239.3480243885944
$ws.Range("G46").Value = 239.3480243885944
$ws.Range("H46").Value = 239.3480243885944
$ws.Range("Q46").Value = 307.8416045893196
$ws.Range("R46").Value = 307.8416045893196
$ws.Range("S46").Value = 307.8416045893196
$ws.Range("T46").Value = 239.3480243885944
$ws.Range("U46").Value = 239.3480243885944
$ws.Range("V46").Value = 239.3480243885944
$ws.Range("W46").Value = 239.3480243885944
$ws.Range("X46").Value = 239.3480243885944
$ws.Range("Y46").Value = 239.3480243885944
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("I2").Value = 44.29520053209296
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("O2").Value = 567.8446588919046
$ws.Range("P2").Value = 870.5779326741233
$ws.Range("Q2").Value = 0
$ws.Range("I5").Value = 44.29520053209296
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 716.8751175230994
$ws.Range("L5").Value = 721.547474042929
$ws.Range("Q5").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 716.8751175230994
$ws.Range("L8").Value = 645.1956588512421
$ws.Range("Q8").Value = 0
$ws.Range("I11").Value = 120.6470157237802
$ws.Range("K11").Value = 716.8751175230994
$ws.Range("L11").Value = 29.34358558166878
$ws.Range("R11").Value = 294.54111633436
$ws.Range("I14").Value = 44.29520053209296
$ws.Range("K14").Value = 496.2158265843866
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 1060.271045550332
$ws.Range("R14").Value = 176.4768357656695
$ws.Range("N17").Value = 1096.663422488788
$ws.Range("O17").Value = 341.7591690772404
$ws.Range("P17").Value = 0
$ws.Range("I20").Value = 120.6470157237802
$ws.Range("J20").Value = 430.3047365861567
$ws.Range("K20").Value = 315.9139665186131
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("Q20").Value = 615.8520732695737
$ws.Range("I23").Value = 44.29520053209296
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 822.5705182964547
$ws.Range("Q23").Value = 615.8520732695737
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 265.4073538855533
$ws.Range("N29").Value = 1096.663422488788
$ws.Range("I32").Value = 120.6470157237802
$ws.Range("K32").Value = 301.7997308240091
$ws.Range("M32").Value = 1060.271045550332
$ws.Range("N32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("M35").Value = 265.4073538855533
$ws.Range("N35").Value = 1096.663422488788
$ws.Range("I38").Value = 120.6470157237802
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 609.557124268908
$ws.Range("M38").Value = 0
$ws.Range("P38").Value = 870.5779326741233
$ws.Range("R38").Value = 176.4768357656695
$ws.Range("I44").Value = 120.6470157237802
$ws.Range("J44").Value = 430.3047365861567
$ws.Range("K44").Value = 315.9139665186112
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("Q44").Value = 615.8520732695737
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("F11").Value = 1.388344292978218
$ws.Range("Y11").Value = 0
$ws.Range("R12").Value = 0
$ws.Range("S12").Value = 1.533699957572765
$ws.Range("T13").Value = 33.55489271692969
$ws.Range("W13").Value = 226.3728098387097
$ws.Range("G14").Value = 1.388344292978218
$ws.Range("T14").Value = 0
$ws.Range("R15").Value = 1.53369995757339
$ws.Range("S15").Value = 0
$ws.Range("F16").Value = 274.3828559677419
$ws.Range("H16").Value = 40.13443591884794
$ws.Range("D19").Value = 285.5362180555555
$ws.Range("H19").Value = 40.13443591884794
$ws.Range("H20").Value = 1.388344292976399
$ws.Range("S20").Value = 0
$ws.Range("H22").Value = 205.3472244622743
$ws.Range("T22").Value = 33.55489271692969
$ws.Range("R24").Value = 1.53369995757339
$ws.Range("S24").Value = 0
$ws.Range("Q25").Value = 505.228266425598
$ws.Range("T25").Value = 33.55489271692969
$ws.Range("H26").Value = 0
$ws.Range("X26").Value = 1.38834429297799
$ws.Range("E28").Value = 280.9809048369565
$ws.Range("I28").Value = 97.40414414470841
$ws.Range("V28").Value = 33.9575216727898
$ws.Range("H29").Value = 0
$ws.Range("V29").Value = 1.388344292978104
$ws.Range("H31").Value = 40.13443591884794
$ws.Range("I31").Value = 97.40414414470841
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = 0
$ws.Range("S32").Value = 1.388344292978275
$ws.Range("H34").Value = 40.13443591884794
$ws.Range("C37").Value = 272.7252466480447
$ws.Range("H37").Value = 40.13443591884794
$ws.Range("G38").Value = 1.388344292978218
$ws.Range("Y38").Value = 0
$ws.Range("R39").Value = 1.53369995757339
$ws.Range("W39").Value = 0
$ws.Range("V40").Value = 33.9575216727898
$ws.Range("X40").Value = 247.4436454301076
$ws.Range("C41").Value = 0
$ws.Range("Y41").Value = 1.388344292978161
$ws.Range("R42").Value = 0
$ws.Range("S42").Value = 1.533699957572765
$ws.Range("G44").Value = 0
$ws.Range("Y44").Value = 1.388344292978161
$ws.Range("I46").Value = 0
$ws.Range("Q46").Value = 505.228266425598
$ws.Range("T46").Value = 130.9590368616381
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B2").Value = 1316006.940470103
$ws.Range("B3").Value = 2606868.406811011
$ws.Range("B4").Value = 3897729.873151923
$ws.Range("B5").Value = 5053973.894141731
$ws.Range("B6").Value = 6210217.915131541
$ws.Range("B7").Value = 7366461.936121343
$ws.Range("B8").Value = 8522705.957111143
$ws.Range("B9").Value = 9678949.978100942
$ws.Range("B10").Value = 10835193.99909074
$ws.Range("B11").Value = 11991438.02008054
$ws.Range("B12").Value = 13147682.04107034
$ws.Range("B13").Value = 14303926.06206015
$ws.Range("B14").Value = 15460170.08304997
$ws.Range("B15").Value = 16616414.10403979
$ws.Range("B16").Value = 17772658.12502962
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 913048.3542411315
$ws.Range("C2").Value = 913048.3542411315
$ws.Range("F2").Value = 817831.1367976689
$ws.Range("H2").Value = 817831.1367976691
$ws.Range("I2").Value = 817831.1367976691
$ws.Range("K2").Value = 817831.136797669
$ws.Range("L2").Value = 817831.1367976689
$ws.Range("M2").Value = 817831.136797669
$ws.Range("O2").Value = 817831.1367976686
$ws.Range("P2").Value = 817831.1367976688
$ws.Range("E4").Value = 26356.19966309947
$ws.Range("G4").Value = 26356.19966309946
$ws.Range("H4").Value = 26356.19966309948
$ws.Range("I4").Value = 26356.19966309947
$ws.Range("L4").Value = 26356.19966309947
$ws.Range("N4").Value = 26356.19966309946
$ws.Range("P4").Value = 26356.19966309947
$ws.Range("C6").Value = 639094.764409381
$ws.Range("F6").Value = 652495.5371345694
$ws.Range("H6").Value = 652495.5371345696
$ws.Range("I6").Value = 652495.5371345696
$ws.Range("K6").Value = 652495.5371345696
$ws.Range("L6").Value = 652495.5371345694
$ws.Range("M6").Value = 652495.5371345696
$ws.Range("O6").Value = 652495.5371345691
$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Range("H2").Value = 0
$ws.Range("S2").Value = 1.388344292978275
$ws.Range("G3").Value = 0
$ws.Range("I3").Value = 1.533699957573035
$ws.Range("I4").Value = 97.40414414470841
$ws.Range("M4").Value = 258.8126637769763
$ws.Range("S5").Value = 0
$ws.Range("X5").Value = 1.38834429297799
$ws.Range("G6").Value = 0
$ws.Range("I6").Value = 1.533699957573035
$ws.Range("O7").Value = 268.0582198167198
$ws.Range("S7").Value = 400
$ws.Range("T7").Value = 235.1986680429646
$ws.Range("H8").Value = 0
$ws.Range("S8").Value = 1.388344292978275
$ws.Range("F9").Value = 1.533699957572935
$ws.Range("I9").Value = 0
$ws.Range("K10").Value = 151.1650359809837
$ws.Range("N10").Value = 168.0209277307011
$ws = $wb.Worksheets.Item("Battery Input ")
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 463.4516620966183
$ws.Range("K2").Value = 747.1142864321608
$ws.Range("L2").Value = 926.8613366834174
$ws.Range("O2").Value = 1487.193025498716
$ws.Range("P2").Value = 1714.888365095753
$ws.Range("Q2").Value = 461.4348583860824
$ws.Range("I4").Value = 0
$ws.Range("M4").Value = 151.1650359809837
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 463.4516620966183
$ws.Range("K5").Value = 1463.989403955261
$ws.Range("L5").Value = 1648.408810726346
$ws.Range("Q5").Value = 461.4348583860824
$ws.Range("O7").Value = 0
$ws.Range("S7").Value = 83.54024160861997
$ws.Range("T7").Value = 36.43098678260854
$ws.Range("J8").Value = 463.4516620966183
$ws.Range("K8").Value = 1463.98940395526
$ws.Range("L8").Value = 1572.05699553466
$ws.Range("Q8").Value = 461.4348583860824
$ws.Range("K10").Value = 168.7357375896766
$ws.Range("N10").Value = 0
$ws.Range("I11").Value = 76.35181519168728
$ws.Range("K11").Value = 1463.98940395526
$ws.Range("L11").Value = 956.2049222650858
$ws.Range("R11").Value = 118.0642805686905
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 1243.330113016547
$ws.Range("L14").Value = 926.8613366834172
$ws.Range("M14").Value = 1547.320742028618
$ws.Range("R14").Value = 0
$ws.Range("N17").Value = 1667.413992057144
$ws.Range("O17").Value = 1261.107535684052
$ws.Range("P17").Value = 844.3104324216301
$ws.Range("I20").Value = 76.35181519168728
$ws.Range("J20").Value = 893.7563986827749
$ws.Range("K20").Value = 1063.028252950774
$ws.Range("L20").Value = 926.861336683417
$ws.Range("M20").Value = 487.0496964782857
$ws.Range("N20").Value = 570.7505695683558
$ws.Range("Q20").Value = 1077.286931655656
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 747.1142864321608
$ws.Range("L23").Value = 926.8613366834174
$ws.Range("M23").Value = 1309.62021477474
$ws.Range("Q23").Value = 1077.286931655656
$ws.Range("J29").Value = 463.4516620966183
$ws.Range("K29").Value = 747.114286432161
$ws.Range("L29").Value = 926.861336683417
$ws.Range("M29").Value = 752.457050363839
$ws.Range("N29").Value = 1667.413992057144
$ws.Range("I32").Value = 76.35181519168728
$ws.Range("K32").Value = 1048.91401725617
$ws.Range("M32").Value = 1547.320742028618
$ws.Range("N32").Value = 570.7505695683558
$ws.Range("P32").Value = 844.3104324216301
$ws.Range("J35").Value = 463.4516620966183
$ws.Range("M35").Value = 752.457050363839
$ws.Range("N35").Value = 1667.413992057144
$ws.Range("I38").Value = 76.35181519168728
$ws.Range("J38").Value = 463.4516620966183
$ws.Range("K38").Value = 1356.671410701069
$ws.Range("M38").Value = 487.0496964782857
$ws.Range("P38").Value = 1714.888365095753
$ws.Range("R38").Value = 0
$ws.Range("I44").Value = 76.35181519168728
$ws.Range("J44").Value = 893.7563986827749
$ws.Range("K44").Value = 1063.028252950772
$ws.Range("M44").Value = 487.0496964782857
$ws.Range("N44").Value = 570.7505695683558
$ws.Range("Q44").Value = 1077.286931655656

Write-Host "Applied changes"